$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style of the existing header cell (H1) onto the new header cells
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Header row values for the two new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for rows 2-30
$data = @{
    2  = @(4, 4)
    3  = @(5, 5)
    4  = @(5, 5)
    5  = @(4, 5)
    6  = @(6, 8)
    7  = @(7, 8)
    8  = @(6, 6)
    9  = @(8, 8)
    10 = @(7, 7)
    11 = @(6, 6)
    12 = @(7, 7)
    13 = @(9, 9)
    14 = @(6, 6)
    15 = @(5, 6)
    16 = @(7, 7)
    17 = @(6, 7)
    18 = @(5, 6)
    19 = @(8, 9)
    20 = @(6, 6)
    21 = @(7, 7)
    22 = @(7, 7)
    23 = @(6, 7)
    24 = @(7, 7)
    25 = @(6, 7)
    26 = @(6, 8)
    27 = @(6, 6)
    28 = @(8, 8)
    29 = @(8, 8)
    30 = @(5, 5)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
